$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the floating point precision of the existing A14 timestamp value
$ws.Cells.Item(14, 1).Value = 44327.78192163195

# Append the new data row (row 15)
$ws.Cells.Item(15, 1).Value = 44328.77718483692
$ws.Cells.Item(15, 2).Value = 74514
$ws.Cells.Item(15, 3).Value = 62614
$ws.Cells.Item(15, 4).Value = 3226
$ws.Cells.Item(15, 5).Value = 2082
$ws.Cells.Item(15, 6).Value = 1476
$ws.Cells.Item(15, 7).Value = 19336
$ws.Cells.Item(15, 8).Value = 1324
$ws.Cells.Item(15, 9).Value = 863
$ws.Cells.Item(15, 10).Value = 219

# The date column (A) uses style index 2 (custom date/time number format).
# Copy the style from the cell above so the new row matches formatting.
$ws.Cells.Item(14, 1).Copy() | Out-Null
$ws.Cells.Item(15, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
